$d = $word.ActiveDocument

# The document contains a single complex field " m:enduserdoc " built from
# fldChar(begin) / instrText / fldChar(end) runs. Replace it with a single
# run containing the literal (escaped) text "{m:enduserdoc}" so the parser
# can pick it up as plain template text instead of a Word field.

$f = $d.Fields.Item(1)
$fieldStart = $f.Code.Start

# Locate the paragraph that currently holds the field.
$paraIndex = 0
$targetIndex = -1
foreach ($para in $d.Paragraphs) {
    $paraIndex = $paraIndex + 1
    if ($fieldStart -ge $para.Range.Start -and $fieldStart -lt $para.Range.End) {
        $targetIndex = $paraIndex
    }
}

# Remove the field (begin/instrText/end runs) from the paragraph, leaving
# an empty paragraph behind.
$f.Delete()

$p = $d.Paragraphs.Item($targetIndex)

# Only target the paragraph's interior (exclude the trailing paragraph
# mark) so the <w:p> element itself - and its rsid attributes - survive;
# InsertXML replaces the contents of exactly the range it's called on.
$interior = $d.Range($p.Range.Start, $p.Range.End - 1)

# Insert a plain run with xml:space="preserve" holding the literal text
# in place of the now-empty field paragraph. InsertXML expects a full
# (flat OPC) WordprocessingML package fragment; only the single <w:r> we
# provide actually lands in the target range.
$flatOpc = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">{m:enduserdoc}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$interior.InsertXML($flatOpc)
